$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.206835666666667
$ws.Range("H2").Value = 9.620507
$ws.Range("I2").Value = 0.1914891004057404
$ws.Range("J2").Value = 0.1914891004057404
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.219409666666667
$ws.Range("N2").Value = 9.658229
$ws.Range("O2").Value = 0.2141077524809907
$ws.Range("P2").Value = 0.2141077524809907
$ws.Range("Q2").Value = 10.32411774467811
$ws.Range("R2").Value = 92.917059702103
$ws.Range("S2").Value = 0.04099930091247985
$ws.Range("T2").Value = 0.04099930091247984

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.206835666666667
$ws.Range("H3").Value = 9.620507
$ws.Range("I3").Value = 0.1914891004057404
$ws.Range("J3").Value = 0.1914891004057404
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.994526666666666
$ws.Range("N3").Value = 14.98358
$ws.Range("O3").Value = 0.3321624117546936
$ws.Range("P3").Value = 0.3321624117546935
$ws.Range("Q3").Value = 16.01662625278444
$ws.Range("R3").Value = 144.14963627506
$ws.Range("S3").Value = 0.06360548141550741
$ws.Range("T3").Value = 0.06360548141550741

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.206835666666667
$ws.Range("H4").Value = 9.620507
$ws.Range("I4").Value = 0.1914891004057404
$ws.Range("J4").Value = 0.1914891004057404
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.295414
$ws.Range("N4").Value = 9.886241999999999
$ws.Range("O4").Value = 0.2191624422141134
$ws.Range("P4").Value = 0.2191624422141134
$ws.Range("Q4").Value = 10.56785115163266
$ws.Range("R4").Value = 95.110660364694
$ws.Range("S4").Value = 0.04196721890230564
$ws.Range("T4").Value = 0.04196721890230564

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.206835666666667
$ws.Range("H5").Value = 9.620507
$ws.Range("I5").Value = 0.1914891004057404
$ws.Range("J5").Value = 0.1914891004057404
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.527049
$ws.Range("N5").Value = 10.581147
$ws.Range("O5").Value = 0.2345673935502023
$ws.Range("P5").Value = 0.2345673935502023
$ws.Range("Q5").Value = 11.310666531281
$ws.Range("R5").Value = 101.795998781529
$ws.Range("S5").Value = 0.04491709917544752
$ws.Range("T5").Value = 0.04491709917544752

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.680435666666667
$ws.Range("H6").Value = 17.041307
$ws.Range("I6").Value = 0.3391946544156194
$ws.Range("J6").Value = 0.3391946544156194
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.219409666666667
$ws.Range("N6").Value = 9.658229
$ws.Range("O6").Value = 0.2141077524809907
$ws.Range("P6").Value = 0.2141077524809907
$ws.Range("Q6").Value = 18.28764949614478
$ws.Range("R6").Value = 164.588845465303
$ws.Range("S6").Value = 0.07262420511049462
$ws.Range("T6").Value = 0.07262420511049461

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.680435666666667
$ws.Range("H7").Value = 17.041307
$ws.Range("I7").Value = 0.3391946544156194
$ws.Range("J7").Value = 0.3391946544156194
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.994526666666666
$ws.Range("N7").Value = 14.98358
$ws.Range("O7").Value = 0.3321624117546936
$ws.Range("P7").Value = 0.3321624117546935
$ws.Range("Q7").Value = 28.37108741545111
$ws.Range("R7").Value = 255.33978673906
$ws.Range("S7").Value = 0.1126677144649919
$ws.Range("T7").Value = 0.1126677144649919

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.680435666666667
$ws.Range("H8").Value = 17.041307
$ws.Range("I8").Value = 0.3391946544156194
$ws.Range("J8").Value = 0.3391946544156194
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.295414
$ws.Range("N8").Value = 9.886241999999999
$ws.Range("O8").Value = 0.2191624422141134
$ws.Range("P8").Value = 0.2191624422141134
$ws.Range("Q8").Value = 18.71938722203267
$ws.Range("R8").Value = 168.474484998294
$ws.Range("S8").Value = 0.07433872884769935
$ws.Range("T8").Value = 0.07433872884769933

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.680435666666667
$ws.Range("H9").Value = 17.041307
$ws.Range("I9").Value = 0.3391946544156194
$ws.Range("J9").Value = 0.3391946544156194
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.527049
$ws.Range("N9").Value = 10.581147
$ws.Range("O9").Value = 0.2345673935502023
$ws.Range("P9").Value = 0.2345673935502023
$ws.Range("Q9").Value = 20.035174937681
$ws.Range("R9").Value = 180.316574439129
$ws.Range("S9").Value = 0.07956400599243346
$ws.Range("T9").Value = 0.07956400599243346

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.798321333333333
$ws.Range("H10").Value = 11.394964
$ws.Range("I10").Value = 0.2268083590101642
$ws.Range("J10").Value = 0.2268083590101642
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.219409666666667
$ws.Range("N10").Value = 9.658229
$ws.Range("O10").Value = 0.2141077524809907
$ws.Range("P10").Value = 0.2141077524809907
$ws.Range("Q10").Value = 12.22835241763956
$ws.Range("R10").Value = 110.055171758756
$ws.Range("S10").Value = 0.04856142799156792
$ws.Range("T10").Value = 0.04856142799156791

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.798321333333333
$ws.Range("H11").Value = 11.394964
$ws.Range("I11").Value = 0.2268083590101642
$ws.Range("J11").Value = 0.2268083590101642
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.994526666666666
$ws.Range("N11").Value = 14.98358
$ws.Range("O11").Value = 0.3321624117546936
$ws.Range("P11").Value = 0.3321624117546935
$ws.Range("Q11").Value = 18.97081718790222
$ws.Range("R11").Value = 170.73735469112
$ws.Range("S11").Value = 0.07533721153494051
$ws.Range("T11").Value = 0.07533721153494051

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.798321333333333
$ws.Range("H12").Value = 11.394964
$ws.Range("I12").Value = 0.2268083590101642
$ws.Range("J12").Value = 0.2268083590101642
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.295414
$ws.Range("N12").Value = 9.886241999999999
$ws.Range("O12").Value = 0.2191624422141134
$ws.Range("P12").Value = 0.2191624422141134
$ws.Range("Q12").Value = 12.51704129836533
$ws.Range("R12").Value = 112.653371685288
$ws.Range("S12").Value = 0.04970787387524299
$ws.Range("T12").Value = 0.04970787387524298

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.798321333333333
$ws.Range("H13").Value = 11.394964
$ws.Range("I13").Value = 0.2268083590101642
$ws.Range("J13").Value = 0.2268083590101642
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.527049
$ws.Range("N13").Value = 10.581147
$ws.Range("O13").Value = 0.2345673935502023
$ws.Range("P13").Value = 0.2345673935502023
$ws.Range("Q13").Value = 13.396865460412
$ws.Range("R13").Value = 120.571789143708
$ws.Range("S13").Value = 0.05320184560841276
$ws.Range("T13").Value = 0.05320184560841276

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.061238666666667
$ws.Range("H14").Value = 12.183716
$ws.Range("I14").Value = 0.242507886168476
$ws.Range("J14").Value = 0.242507886168476
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.219409666666667
$ws.Range("N14").Value = 9.658229
$ws.Range("O14").Value = 0.2141077524809907
$ws.Range("P14").Value = 0.2141077524809907
$ws.Range("Q14").Value = 13.07479102210711
$ws.Range("R14").Value = 117.673119198964
$ws.Range("S14").Value = 0.05192281846644833
$ws.Range("T14").Value = 0.05192281846644832

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.061238666666667
$ws.Range("H15").Value = 12.183716
$ws.Range("I15").Value = 0.242507886168476
$ws.Range("J15").Value = 0.242507886168476
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.994526666666666
$ws.Range("N15").Value = 14.98358
$ws.Range("O15").Value = 0.3321624117546936
$ws.Range("P15").Value = 0.3321624117546935
$ws.Range("Q15").Value = 20.28396482036445
$ws.Range("R15").Value = 182.55568338328
$ws.Range("S15").Value = 0.08055200433925368
$ws.Range("T15").Value = 0.08055200433925366

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.061238666666667
$ws.Range("H16").Value = 12.183716
$ws.Range("I16").Value = 0.242507886168476
$ws.Range("J16").Value = 0.242507886168476
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 3.295414
$ws.Range("N16").Value = 9.886241999999999
$ws.Range("O16").Value = 0.2191624422141134
$ws.Range("P16").Value = 0.2191624422141134
$ws.Range("Q16").Value = 13.38346275947467
$ws.Range("R16").Value = 120.451164835272
$ws.Range("S16").Value = 0.0531486205888654
$ws.Range("T16").Value = 0.05314862058886539

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.061238666666667
$ws.Range("H17").Value = 12.183716
$ws.Range("I17").Value = 0.242507886168476
$ws.Range("J17").Value = 0.242507886168476
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 3.527049
$ws.Range("N17").Value = 10.581147
$ws.Range("O17").Value = 0.2345673935502023
$ws.Range("P17").Value = 0.2345673935502023
$ws.Range("Q17").Value = 14.324187778028
$ws.Range("R17").Value = 128.917690002252
$ws.Range("S17").Value = 0.05688444277390857
$ws.Range("T17").Value = 0.05688444277390856
